$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D32").Value = "Feature importance (in Decision Tree, RF)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/363"

$ws.Range("D51").Value = "[python] UnicodeDecodeError: 'cp949' codec can't decode byte 0xed in position 135: illegal multibyte sequence 에러 해결법"
$ws.Range("E51").Value = "https://bskyvision.com/1262"
